$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "La ridicula idea de no volverte a ver"
$ws.Range("B12").Value = "Rosa Montero"
# Editorial (C) is left blank for this entry, same as other rows without a
# listed publisher (e.g. C8:C11) - touch the cell so it is materialized as
# part of row 12 without pulling in any new formatting.
$ws.Range("C12").Style = "Normal"
